$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("10G108024", "Bag Poly - Turkey 10x8x24 (1mil)", "2", "33.98", "67.96"),
    @("K8", "Wrap Poly 8x10.75", "2", "57.14", "114.28"),
    @("3M84CC", "Scrubbies - Steel", "1", "48.53", "48.53"),
    @("86N", "Scrubbies - Green", "1", "6.56", "6.56"),
    @("765004", "Extreme - AllTemp Machine Detergent", "1", "131.29", "131.29"),
    @("764804", "48 SparClean Dish Manual", "1", "102.51", "102.51")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($row, $col)
        # Force text storage (source cells are inline/shared strings, even
        # for numeric-looking values like "2" or "33.98"), then restore the
        # default "Normal" style so no stray style index is left on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$j]
        $cell.Style = "Normal"
    }
}
